# Widen the two columns of the "Data documentation" table to make room
# for the new flights row(s) added to the plotly viz write-up.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Widths are expressed in points; the document stores widths in twips
# (1 pt = 20 twips), so 3384/20 = 169.2 and 5627/20 = 281.35.
$table.Columns.Item(1).Width = 169.2
$table.Columns.Item(2).Width = 281.35
